$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge E14:E16 so each row gets its own value/box-border style ---
$ws.Range("E14:E16").UnMerge()

# --- Row 14: Part of the "24-8-23 & 25-8-23" work entry (top of the E14:E16 box) ---
$ws.Range("E14").Value = 2
$ws.Range("E14").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("E14").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("E14").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws.Range("E14").HorizontalAlignment = -4108     # xlCenter
$ws.Range("E14").VerticalAlignment = -4108       # xlCenter

$ws.Range("F14").Value = 1.5
$ws.Range("G14").Value = "24-8-23 & 25-8-23"

# --- Row 15 (middle of the E14:E16 box) ---
$ws.Range("E15").Value = 2
$ws.Range("E15").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("E15").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("E15").HorizontalAlignment = -4108
$ws.Range("E15").VerticalAlignment = -4108

# --- Row 16 (bottom of the E14:E16 box) ---
$ws.Range("E16").Value = 2
$ws.Range("E16").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("E16").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("E16").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").VerticalAlignment = -4108

# --- Daily work-done hours added for rows 21-26 ---
$ws.Range("E21").Value = 2
$ws.Range("E22").Value = 2
$ws.Range("E23").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("E26").Value = 2

# --- Widen column G so the new combined date text fits ---
$ws.Columns.Item(7).ColumnWidth = 17

# --- Move the active-cell cursor, mirroring the author's last-edited cell ---
$ws.Range("E25").Select()
